$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.191.35"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.055.03"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "248.63"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").Value = "57.81"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "15.92"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "0.874"
$ws.Range("E13").Value = "  +6.17%  "
$ws.Range("D14").Value = "2.352.77"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "2.048.24"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "18.21"
$ws.Range("E17").Value = "  +15.22%  "
$ws.Range("D18").Value = "37.229.37"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "74.92"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "237.39"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -4.37%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.49"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "169.33"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "20.11"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D30").Value = "4.83"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "0.0890"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  +15.26%  "
$ws.Range("E40").Value = "  +11.60%  "
$ws.Range("D41").Value = "0.0981"
$ws.Range("E41").Value = "  -15.44%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "17.32"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "1.14"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").Value = "96.11"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "1.269.00"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "2.239.40"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "44.01"
$ws.Range("E51").Value = "  -0.70%  "
